# Fix load_words function to read from words.xlsx and add new words.xlsx file
# -- adds three new verb-conjugation rows (落とす/落ちる/残る groups) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order across the sheet is:
# A = Dictionary Form, B = Te Form, C = Ta Form, D = Nai Form,
# E = Masu Form, F = Volitional Form, G = Potential Form

$newRows = @(
    @("落とす",   "落として",   "落とした",   "落とさない",   "落とします",   "落とそう",   "落とせる"),
    @("落ちる",   "落ちて",     "落ちた",     "落ちない",     "落ちます",     "落ちよう",   "落ちられる"),
    @("残る",     "残って",     "残った",     "残らない",     "残ります",     "残ろう",     "残れる")
)

$startRow = 104
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]

    $ws.Rows.Item($r).RowHeight = 18.75

    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        $cell.Value = $rowValues[$c]
        # Match the workbook's existing convention for Japanese-text cells
        # (style index 3 / Yu Gothic font) used throughout the table, except
        # the very first Dictionary-Form cell of this new block, which stays
        # on the default style just like row 101's block start.
        if (-not ($r -eq $startRow -and $c -eq 0)) {
            $cell.Font.Name = "Yu Gothic"
        }
    }
}

# Reflect the user's final scroll/selection position from the edit session.
$ws.Range("F107").Select()

Write-Host "Added verb rows 104-106 (落とす/落ちる/残る)"
